# Update cryptos list: refresh Price (D) and Volume(1h) (E) columns
# for rows 2-51 on the active sheet, matching the scraped snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.154.75'
$ws.Range('E2').Value = '  +1.94%  '
$ws.Range('D3').Value = '3.776.77'
$ws.Range('E3').Value = '  -0.24%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '626.09'
$ws.Range('E5').Value = '  +3.90%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '165.74'
$ws.Range('E6').Value = '  +1.64%  '
$ws.Range('D7').Value = '3.774.11'
$ws.Range('E7').Value = '  -0.21%  '
$ws.Range('E8').Value = '  -0.08%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.522'
$ws.Range('E9').Value = '  +1.51%  '
$ws.Range('E10').Value = '  +1.21%  '
$ws.Range('E11').Value = '  +3.07%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '6.76'
$ws.Range('E12').Value = '  -0.80%  '
$ws.Range('E13').Value = '  +0.45%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '35.67'
$ws.Range('E14').Value = '  +1.78%  '
$ws.Range('D15').Value = '4.406.01'
$ws.Range('E15').Value = '  -0.32%  '
$ws.Range('D16').Value = '3.843.02'
$ws.Range('E16').Value = '  +1.45%  '
$ws.Range('D17').Value = '69.140.95'
$ws.Range('E17').Value = '  +1.92%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '17.65'
$ws.Range('E18').Value = '  -2.84%  '
$ws.Range('E19').Value = '  -1.15%  '
$ws.Range('E20').Value = '  +0.55%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '467.40'
$ws.Range('E21').Value = '  +2.04%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.58'
$ws.Range('E22').Value = '  +1.59%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.707'
$ws.Range('E23').Value = '  +2.37%  '
$ws.Range('E24').Value = '  +2.91%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '83.13'
$ws.Range('E25').Value = '  +0.09%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.05'
$ws.Range('E26').Value = '  +1.72%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.16'
$ws.Range('E27').Value = '  +3.74%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.04'
$ws.Range('E28').Value = '  +1.56%  '
$ws.Range('E29').Value = '  -0.08%  '
$ws.Range('D30').Value = '3.922.90'
$ws.Range('E30').Value = '  -0.32%  '
$ws.Range('E31').Value = '  +2.65%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.24'
$ws.Range('E32').Value = '  +2.46%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '7.16'
$ws.Range('E33').Value = '  -0.68%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '28.81'
$ws.Range('E34').Value = '  -0.45%  '
$ws.Range('E35').Value = '  +20.88%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.999'
$ws.Range('E36').Value = '  +0.16%  '
$ws.Range('D37').Value = '3.724.59'
$ws.Range('E37').Value = '  -0.26%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '8.96'
$ws.Range('E38').Value = '  +0.49%  '
$ws.Range('E39').Value = '  +2.40%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.36'
$ws.Range('E40').Value = '  +5.48%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.83'
$ws.Range('E41').Value = '  +0.60%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.969'
$ws.Range('E42').Value = '  -0.56%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.00'
$ws.Range('E43').Value = '  +0.00%  '
$ws.Range('E44').Value = '  -0.07%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '153.75'
$ws.Range('E45').Value = '  +1.34%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '43.31'
$ws.Range('E46').Value = '  -0.86%  '
$ws.Range('E47').Value = '  +0.92%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '46.66'
$ws.Range('E48').Value = '  -0.92%  '
$ws.Range('E49').Value = '  +4.01%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.41'
$ws.Range('E50').Value = '  +1.55%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.37'
$ws.Range('E51').Value = '  +0.36%  '
